# Asset addition and validation added
#
# This script reproduces, via Excel COM automation, the changes described by
# the target diff:
#   1. Manufacture sheet: remove the stray row 17 ("S") that was left over
#      from earlier test data (also garbage-collects the "S" shared string).
#   2. Location sheet: append four new location rows (LID03, LID11, LID12,
#      and a duplicate LID010 row).
#   3. User sheet: append a new "Technician" user row (Auto-UID003) with its
#      own mailto hyperlink.
#   4. Add a brand-new "Asset" worksheet (after "User") with a header row and
#      one sample asset row, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Manufacture: delete the orphan row 17 (C17 = "S")
# ---------------------------------------------------------------------------
$wsManufacture = $wb.Worksheets.Item("Manufacture")
$wsManufacture.Rows.Item(17).Delete()
$wsManufacture.Activate()
$wsManufacture.Range("D24").Select()

# ---------------------------------------------------------------------------
# 2. Location: append rows 4-6
# ---------------------------------------------------------------------------
$wsLocation = $wb.Worksheets.Item("Location")

$wsLocation.Range("A4").Value = "Auto-LID03"
$wsLocation.Range("B4").Value = "'12.35"
$wsLocation.Range("C4").Value = "'34.569"
$wsLocation.Range("D4").Value = "'05:03:02"
$wsLocation.Range("E4").Value = "Automation"
$wsLocation.Range("F4").Value = "wimate"
$wsLocation.Range("G4").Value = "Auto-LTID02"

$wsLocation.Range("A5").Value = "Auto-LID11"
$wsLocation.Range("B5").Value = "'12.36"
$wsLocation.Range("C5").Value = "'34.570"
$wsLocation.Range("D5").Value = "'05:03:03"
$wsLocation.Range("E5").Value = "Automation"
$wsLocation.Range("F5").Value = "wimate"
$wsLocation.Range("G5").Value = "Auto-LTID01"

$wsLocation.Range("A6").Value = "Auto-LID12"
$wsLocation.Range("B6").Value = "'12.37123"
$wsLocation.Range("C6").Value = "'34.571"
$wsLocation.Range("D6").Value = "'05:03:04"
$wsLocation.Range("E6").Value = "Automation"
$wsLocation.Range("F6").Value = "wimate"
$wsLocation.Range("G6").Value = "Auto-LTID01"

$wsLocation.Activate()
$wsLocation.Range("H19").Select()

# ---------------------------------------------------------------------------
# 3. User: append row 4 (new Technician user) + hyperlink
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")

$wsUser.Range("A4").Value = "Auto-UID003"
$wsUser.Range("B4").Value = "Disable"
$wsUser.Range("C4").Value = "Disable"
$wsUser.Range("D4").Value = "Automation"
$wsUser.Range("E4").Value = "user4untangled@gmail.com"
$wsUser.Range("F4").Value = "UID_001"
$wsUser.Range("G4").Value = "Disable"
$wsUser.Range("H4").Value = "Disable"
$wsUser.Range("I4").Value = "'+919483460652"
$wsUser.Range("J4").Value = "Invisible"
$wsUser.Range("K4").Value = "Automation"
$wsUser.Range("L4").Value = "Auto_tester3"
$wsUser.Range("M4").Value = "Auto_002"
$wsUser.Range("N4").Value = "Auto-125"
$wsUser.Range("O4").Value = "Auto-SD124"
$wsUser.Range("P4").Value = "Role2"
$wsUser.Range("Q4").Value = "untangleds"
$wsUser.Range("R4").Value = "Technician"

$wsUser.Hyperlinks.Add($wsUser.Range("E4"), "mailto:user4untangled@gmail.com") | Out-Null

$wsUser.Activate()
$wsUser.Range("R6").Select()

# ---------------------------------------------------------------------------
# 4. New "Asset" sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAsset = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsAsset.Name = "Asset"

$wsAsset.Range("A1").Value = "Asset ID"
$wsAsset.Range("B1").Value = "Tag ID"
$wsAsset.Range("C1").Value = "Asset type id"
$wsAsset.Range("D1").Value = "department ID"
$wsAsset.Range("E1").Value = "Sub department ID"
$wsAsset.Range("F1").Value = "AMC Expiry"
$wsAsset.Range("G1").Value = "bought price"
$wsAsset.Range("H1").Value = "current price"
$wsAsset.Range("I1").Value = "location ID"
$wsAsset.Range("J1").Value = "Vendor Company ID"
$wsAsset.Range("K1").Value = "Asset Type"
$wsAsset.Range("L1").Value = "Manufacturer ID"
$wsAsset.Range("M1").Value = "Model ID"
$wsAsset.Range("N1").Value = "Meta Data"
$wsAsset.Range("O1").Value = "Installation Date"
$wsAsset.Range("P1").Value = "Latest service Date"
$wsAsset.Range("Q1").Value = "Next service Date"
$wsAsset.Range("R1").Value = "group"
$wsAsset.Range("S1").Value = "Technician ID"
$wsAsset.Range("T1").Value = "User Notify ID"

$wsAsset.Range("A2").Value = "Auto-Asset"
$wsAsset.Range("B2").Value = "Auto-TID01"
$wsAsset.Range("C2").Value = "Auto-ATID01"
$wsAsset.Range("D2").Value = "Auto-123"
$wsAsset.Range("E2").Value = "Auto-SD131"

$wsAsset.Range("F2").Value = "'01/10/2023"
$wsAsset.Range("F2").NumberFormat = "m/d/yyyy"

$wsAsset.Range("G2").Value = "'25000"
$wsAsset.Range("H2").Value = "'20000"

$wsAsset.Range("I2").Value = "Auto-LID11"
$wsAsset.Range("J2").Value = "VC01"
$wsAsset.Range("K2").Value = "Fixed"
$wsAsset.Range("L2").Value = "Auto-MID140"
$wsAsset.Range("M2").Value = "Auto-MDID132"
$wsAsset.Range("N2").Value = "Automation"

$wsAsset.Range("O2").Value = "'1/17/2023"
$wsAsset.Range("P2").Value = "'3/28/2023"
$wsAsset.Range("Q2").Value = "'1/9/2023"

$wsAsset.Range("R2").Value = "Automation"
$wsAsset.Range("S2").Value = "UID_1234"
$wsAsset.Range("T2").Value = "UID_001"

$wsAsset.Activate()
$wsAsset.Range("D17").Select()
